$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Fix typo: "Leased Employed" -> "Leased Employee" in the Assigned To column (D)
$ws.Range("D4").Value = "Oboe Andrea (Leased Employee)"
$ws.Range("D5").Value = "Oboe Andrea (Leased Employee)"
$ws.Range("D6").Value = "Oboe Andrea (Leased Employee)"

# Update selection to match author's last cursor position
$ws.Range("E5").Select()
